$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 13114.9679433723
$ws.Range("E2").Value = 9788.61473519149
$ws.Range("F2").Value = 72.6967782734894
$ws.Range("C3").Value = 9223.07598166684
$ws.Range("E3").Value = 9473.67959664485
$ws.Range("F3").Value = 275.870649096321
$ws.Range("C4").Value = 13939.2556180087
$ws.Range("E4").Value = 9893.19026967497
$ws.Range("F4").Value = 489.857745320151
$ws.Range("C5").Value = 14212.4831382548
$ws.Range("E5").Value = 9966.80923640644
$ws.Range("F5").Value = 504.30968227755
$ws.Range("C6").Value = 12523.2976755524
$ws.Range("E6").Value = 9581.76793397161
$ws.Range("F6").Value = 417.883567063499
$ws.Range("C7").Value = 8606.69930379694
$ws.Range("E7").Value = 8698.44036097227
$ws.Range("F7").Value = 217.886652698717
$ws.Range("C9").Value = 12383.1420338857
$ws.Range("F9").Value = 392.171753298609
$ws.Range("C10").Value = 12337.0860976626
$ws.Range("F10").Value = 390.252755955979
$ws.Range("C11").Value = 12392.6033065818
$ws.Range("F11").Value = 392.565972994278
$ws.Range("C12").Value = 12690.1874016169
$ws.Range("F12").Value = 404.965310287407
$ws.Range("C13").Value = 11801.34359549
$ws.Range("F13").Value = 367.928481674871
$ws.Range("C14").Value = 8194.39458171517
$ws.Range("F14").Value = 200.355151279665
$ws.Range("C15").Value = 7966.17118587628
$ws.Range("F15").Value = 190.328514127228
